# Apply the "additional weekly and monthly reports with filtered" edit:
#   - Remove the ContactNo column (column C) entirely.
#   - Add a "Total:" label in A5 (the row that already holds the 7000 sum in B5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole ContactNo column (C) - shifts nothing left of it, just removes it.
$ws.Range("C1:C4").Delete()

# Label the totals row.
$ws.Range("A5").Value = "Total:"
